$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("DataSet")

# Insert a new row above current row 5 (shifts rows 5..46 down to 6..47).
# Excel copies formatting from the row above (row 4), which already matches
# the style pattern required for the new row (s="3"/s="11"/s="6" on the
# relevant columns).
$ws.Rows.Item(5).Insert()

# Populate the new row 5 with the employee/user credentials.
$ws.Range("A5").Value = "Employeeuser"
$ws.Range("B5").Value = "vpilli@helenoftroy.com"
$ws.Range("C5").Value = "vpilli@helenoftroy.com"
$ws.Range("E5").Value = "Lotuswave123"
$ws.Range("F5").Value = "Lotuswave123"

# Add mailto hyperlinks on the email cells, matching the pattern used for
# the equivalent cells in row 4 (B4/C4 -> mailto:<email>).
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:vpilli@helenoftroy.com")
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:vpilli@helenoftroy.com")

# Adding a hyperlink re-styles the cell with the built-in "Hyperlink" look;
# restore the original (non-hyperlink-colored) style these cells should keep,
# matching row 4's equivalent B4/C4 cells which already carry that style.
$ws.Range("B5").Style = $ws.Range("B4").Style
$ws.Range("C5").Style = $ws.Range("C4").Style

# Update the sheet selection to reflect where the user ended up after the edit.
$ws.Range("I7").Select()
